$d = $word.ActiveDocument

# --- 1. Insert the new "[[IN THE DIGITAL VERSION SCROLLJACK THIS]]" paragraph ---
# It goes right after the existing empty BodyText paragraph that follows
# "The latter is shown below." (and before the "When pediatric endocrinologist
# Charmian Quigley ..." paragraph), keeping the empty paragraph untouched.
$r = $d.Content
$r.Find.Execute("The latter is shown below.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.MoveEnd(1, 1)
$r.Collapse(0)
$r.InsertAfter("`r[[IN THE DIGITAL VERSION SCROLLJACK THIS]]")

# --- 2. Style tweaks ---

# Compact: center-align.
$sCompact = $d.Styles("Compact")
$sCompact.ParagraphFormat.Alignment = 1

# Caption: center-align + grey color (keeps existing italic).
$sCaption = $d.Styles("Caption")
$sCaption.ParagraphFormat.Alignment = 1
$sCaption.Font.Color = 8421504

# Figure: center-align (new pPr).
$sFigure = $d.Styles("Figure")
$sFigure.ParagraphFormat.Alignment = 1

# CaptionChar: Times New Roman, italic, grey color.
$sCaptionChar = $d.Styles("CaptionChar")
$sCaptionChar.Font.Name = "Times New Roman"
$sCaptionChar.Font.Italic = $true
$sCaptionChar.Font.Color = 8421504

# VerbatimChar: add italic + grey color (keeps Consolas + sz 22).
$sVerbatimChar = $d.Styles("VerbatimChar")
$sVerbatimChar.Font.Italic = $true
$sVerbatimChar.Font.Color = 8421504

# FootnoteReference: Times New Roman, italic, grey color (keeps superscript).
$sFootnoteRef = $d.Styles("FootnoteReference")
$sFootnoteRef.Font.Name = "Times New Roman"
$sFootnoteRef.Font.Italic = $true
$sFootnoteRef.Font.Color = 8421504

# Hyperlink: Times New Roman, italic (keeps accent1 color).
$sHyperlink = $d.Styles("Hyperlink")
$sHyperlink.Font.Name = "Times New Roman"
$sHyperlink.Font.Italic = $true
